# #5: cash & deposit done
# Add bank / deposit_type / currency metadata columns (G:M) to the
# "存款" (deposit) sheet, matching the property_category/category/date/
# legislator_name/legislator_id/source_file/index columns already used
# on the other property sheets (stock, land, building, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # 存款 (deposits)

# ---- Row 1 : header labels --------------------------------------------
$ws.Cells.Item(1, 7).Value  = "property_category"
$ws.Cells.Item(1, 8).Value  = "category"
# column I holds a literal "yyyy-mm-dd" string, not a real date value -
# force text storage so it isn't reinterpreted as a date serial number.
$dateHeader = $ws.Cells.Item(1, 9)
$dateHeader.NumberFormat = "@"
$dateHeader.Value2 = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# ---- Row 2 : 合作金庫商業銀行景美分行 ----------------------------------
$ws.Cells.Item(2, 7).Value  = "deposit"
$ws.Cells.Item(2, 8).Value  = "normal"
$date2 = $ws.Cells.Item(2, 9)
$date2.NumberFormat = "@"
$date2.Value2 = "2012-05-01"
$ws.Cells.Item(2, 10).Value = "翁重鈞"
$ws.Cells.Item(2, 11).Value = 551
$ws.Cells.Item(2, 12).Value = "tmp6aad1"
$ws.Cells.Item(2, 13).Value = 82

# ---- Row 3 : 臺灣銀行群賢分行 ------------------------------------------
$ws.Cells.Item(3, 7).Value  = "deposit"
$ws.Cells.Item(3, 8).Value  = "normal"
$date3 = $ws.Cells.Item(3, 9)
$date3.NumberFormat = "@"
$date3.Value2 = "2012-05-01"
$ws.Cells.Item(3, 10).Value = "翁重鈞"
$ws.Cells.Item(3, 11).Value = 551
$ws.Cells.Item(3, 12).Value = "tmp6aad1"
$ws.Cells.Item(3, 13).Value = 83

# ---- formatting: match the look of the existing columns ---------------
# header row (B1) uses the bold/bordered style -> copy onto G1:M1
$ws.Cells.Item(1, 2).Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)   # xlPasteFormats

# data rows (B2/B3) use the plain bordered style -> copy onto G:M
$ws.Cells.Item(2, 2).Copy()
$ws.Range("G2:M2").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(3, 2).Copy()
$ws.Range("G3:M3").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false
